$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update header row: remove "-1700" suffix from the day/time headers
$ws.Range("A1").Value = "Monday 1600"
$ws.Range("B1").Value = "Tuesday 1600"
$ws.Range("C1").Value = "Wednesday 1600"
$ws.Range("D1").Value = "Thursday 1600"
$ws.Range("E1").Value = "Friday 1600"

# Update selection to match saved state
$ws.Range("E1").Select()
